# Apply "Update values for FTLR & Tarot" edits
$wb = $excel.ActiveWorkbook

$wsCode = $wb.Worksheets.Item("SAD-Code")
$wsProj = $wb.Worksheets.Item("Projects")

# --- Projects sheet: raw counts updated ---
$wsProj.Range("C4").Value = 1473
$wsProj.Range("E4").Value = 730

# --- SAD-Code sheet: raw metric inputs updated ---
# Block 1 (rows 4-8)
$wsCode.Range("D4").Value = 0.18
$wsCode.Range("F4").Value = 0.9

# Block 2 (rows 15-19)
$wsCode.Range("F15").Value = 0.9

# Block 3 (rows 26-30)
$wsCode.Range("C26").Value = 0.07
$wsCode.Range("D26").Value = 0.18
$wsCode.Range("F26").Value = 0.91
$wsCode.Range("G26").Value = 0.93

# Block 4 (rows 37-41)
$wsCode.Range("D37").Value = 0.16
$wsCode.Range("H37").Value = 0.05

# Block 5 (rows 48-52)
$wsCode.Range("H48").Value = 0.06

# Recalculate all dependent formulas (shared-formula averages & weighted averages)
$excel.CalculateFull()

# --- Restore view/selection state to match the authored edit ---
$wsProj.Range("H5").Select()

$wsCode.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsCode.Range("I13").Select()

$wb.Save()
